$d = $word.ActiveDocument

# Fix typo: "le nombres de bits" -> "le nombre de bits"
$d.Content.Find.Execute("le nombres de bits", $true, $false, $false, $false, $false,
                         $true, 1, $false, "le nombre de bits", 2)

# Remove the stray "_GoBack" bookmark left over from a previous edit session
try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
} catch {
}
